$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: tuition years value changes 30 -> 15 (literal)
$ws.Range("C5").Value = 15

# Row 6: C6 was "=C5" formula; now a plain literal value 15
$ws.Range("C6").Value = 15

# Row 7: C7 was a literal 30; now a formula referencing C6
$ws.Range("C7").Formula = "=C6"

# Row 9: stipend amount 25000 -> 30000
$ws.Range("B9").Value = 30000

# Row 10: C10 was a literal 30; now a formula referencing C7
$ws.Range("C10").Formula = "=C7"

# Row 11: new "Summer salary" line item
$ws.Range("A11").Value = "Summer salary"
$ws.Range("B11").Formula = "=ROUND(75000/9,0)"
$ws.Range("D11").Value = 5
$ws.Range("E11").Formula = "=B11*D11"

# Update the selected cell to match the author's final cursor position
$ws.Range("A12").Select()
